$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "level_id" column (numeric id, was column A) is being dropped from
# the template. "level_kode" and "level_nama" shift one column to the
# left so they become columns A and B.

$ws.Range("A1").Value = "level_kode"
$ws.Range("B1").Value = "level_nama"

$ws.Range("A2").Value = "ADM"
$ws.Range("B2").Value = "Administrator"

$ws.Range("A3").Value = "MNG"
$ws.Range("B3").Value = "Manager"

$ws.Range("A4").Value = "STF"
$ws.Range("B4").Value = "Staff/Kasir"

# Drop the now-empty former level_id/level_nama column C so the used
# range shrinks back down to A1:B4.
$ws.Range("C1:C4").Clear()

# Match the saved selection.
$ws.Range("C3").Select()
